$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: Average of SW(S*)/SW(OPT)  -- format this first so its xf/font land at the lower index, matching target order
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"
$f14 = $ws.Range("B14").Font
$f14.Bold = $true
$f14.Size = 12
$ws.Range("B14").VerticalAlignment = -4108

# Row 15: Average of SC(S*)/SC(OPT)
$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

# Row 16: Worst of SW(S*)/SW(OPT)
$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

# Row 17: Worst of SC(S*)/SC(OPT)
$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Copy B14's formatting onto B15:B17 so they all share the same style index
$ws.Range("B14").Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)

# Row 12: average of column J (the "k" parameter), bold like header style
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$ws.Range("J12").Font.Bold = $true

$ws.Range("J12").Select()

# Match the page setup attributes added by Excel on save
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
